$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.314.59'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '1.868.81'
$ws.Range("E3").Value = '  -0.66%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.05%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '243.74'
$c.ClearFormats()
$ws.Range("E5").Value = '  -2.28%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.ClearFormats()

$ws.Range("E7").Value = '  -0.47%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2875'
$c.ClearFormats()
$ws.Range("E8").Value = '  -2.08%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06482'
$c.ClearFormats()
$ws.Range("E9").Value = '  -0.75%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '21.58'
$c.ClearFormats()
$ws.Range("E10").Value = '  -1.54%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07787'
$c.ClearFormats()
$ws.Range("E11").Value = '  +0.49%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '96.67'
$c.ClearFormats()
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").Value = '1.864.35'
$ws.Range("E13").Value = '  -0.83%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.7188'
$c.ClearFormats()
$ws.Range("E14").Value = '  -2.64%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.131'
$c.ClearFormats()
$ws.Range("E15").Value = '  -1.94%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '281.57'
$c.ClearFormats()
$ws.Range("E16").Value = '  +2.75%  '

$ws.Range("D17").Value = '30.290.78'
$ws.Range("E17").Value = '  -1.33%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '13.00'
$c.ClearFormats()
$ws.Range("E18").Value = '  -1.22%  '

$ws.Range("E19").Value = '  +0.01%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000007469'
$c.ClearFormats()
$ws.Range("E20").Value = '  -0.75%  '

$ws.Range("D21").Value = '2.108.08'
$ws.Range("E21").Value = '  -0.96%  '

$ws.Range("E22").Value = '  +0.00%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.249'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.61%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.255'
$c.ClearFormats()
$ws.Range("E24").Value = '  +0.36%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '162.23'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.85%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.986'
$c.ClearFormats()
$ws.Range("E26").Value = '  -2.56%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.69'
$c.ClearFormats()
$ws.Range("E27").Value = '  -0.82%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.879'
$c.ClearFormats()
$ws.Range("E28").Value = '  -1.35%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.09639'
$c.ClearFormats()
$ws.Range("E29").Value = '  -0.68%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.310'
$c.ClearFormats()
$ws.Range("E30").Value = '  -2.45%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.210'
$c.ClearFormats()
$ws.Range("E32").Value = '  -1.61%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.124'
$c.ClearFormats()
$ws.Range("E33").Value = '  -0.55%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04783'
$c.ClearFormats()
$ws.Range("E34").Value = '  -1.63%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.118'
$c.ClearFormats()
$ws.Range("E35").Value = '  -0.56%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6845'
$c.ClearFormats()
$ws.Range("E36").Value = '  -1.93%  '

$ws.Range("E37").Value = '  -0.30%  '

$ws.Range("E38").Value = '  -0.75%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.836'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.83%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '75.56'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.94%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.225'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.06%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.928'
$c.ClearFormats()
$ws.Range("E42").Value = '  -4.89%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.4200'
$c.ClearFormats()
$ws.Range("E43").Value = '  -1.01%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.05%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.8241'
$c.ClearFormats()
$ws.Range("E45").Value = '  -1.93%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '100.39'
$c.ClearFormats()
$ws.Range("E46").Value = '  -2.20%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.607'
$c.ClearFormats()
$ws.Range("E47").Value = '  +2.63%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '6.978'
$c.ClearFormats()
$ws.Range("E48").Value = '  -0.87%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '34.96'
$c.ClearFormats()
$ws.Range("E49").Value = '  -1.79%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05759'
$c.ClearFormats()
$ws.Range("E50").Value = '  +0.19%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '883.72'
$c.ClearFormats()
$ws.Range("E51").Value = '  -3.46%  '
